# TutorialsNinjaTestData.xlsx - trim the Login test-data sheet down to a
# single credential row and bump the password value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop every hyperlink up front - only one survives the cleanup below and
# it is rebuilt from scratch once the surviving row is in place.
$ws.Hyperlinks.Delete() | Out-Null

# Row 6 (dheeruvish1612@gmail.com) carried the "last row" look (no top
# border on the e-mail cell, outer border on the password cell). Copy that
# formatting onto row 2 so the single row left behind keeps that look.
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Bump the surviving row's password value.
$ws.Range("B2").Value = 123456

# Remove the now-redundant credential rows (3-6); this also clears out the
# shared strings / relationships tied to the removed e-mail addresses.
$ws.Range("A3:A6").EntireRow.Delete() | Out-Null

# Re-create the single surviving hyperlink (kept pointing at the first
# test e-mail). Adding it restyles the cell, so reapply the "Hyperlink"
# look afterwards to land back on the font/border combo already in use.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:dheeruvish1608@gmail.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

# Leave the saved cursor position where the author left it.
$ws.Range("D7").Select() | Out-Null
